{"js": "// Replace each two-digit-division answer cell's text with its new value.\n// The mapping below was derived from the target diff: each entry is\n// [oldText, newText] for a <w:t> run inside a table cell.\nconst replacements = [\n  [\"42\u00f76=7, 0\", \"40\u00f72=20, 0\"],\n  [\"87\u00f78=10, 7\", \"34\u00f76=5, 4\"],\n  [\"13\u00f77=1, 6\", \"47\u00f73=15, 2\"],\n  [\"60\u00f79=6, 6\", \"53\u00f73=17, 2\"],\n  [\"17\u00f78=2, 1\", \"84\u00f74=21, 0\"],\n  [\"31\u00f72=15, 1\", \"32\u00f74=8, 0\"],\n  [\"29\u00f73=9, 2\", \"85\u00f73=28, 1\"],\n  [\"99\u00f78=12, 3\", \"65\u00f74=16, 1\"],\n  [\"54\u00f78=6, 6\", \"42\u00f73=14, 0\"],\n  [\"31\u00f76=5, 1\", \"46\u00f79=5, 1\"],\n  [\"92\u00f74=23, 0\", \"71\u00f79=7, 8\"],\n  [\"13\u00f79=1, 4\", \"45\u00f78=5, 5\"],\n  [\"83\u00f79=9, 2\", \"94\u00f76=15, 4\"],\n  [\"83\u00f73=27, 2\", \"11\u00f79=1, 2\"],\n  [\"93\u00f76=15, 3\", \"20\u00f75=4, 0\"],\n  [\"90\u00f72=45, 0\", \"69\u00f78=8, 5\"],\n  [\"12\u00f77=1, 5\", \"98\u00f72=49, 0\"],\n  [\"16\u00f77=2, 2\", \"68\u00f77=9, 5\"],\n  [\"24\u00f74=6, 0\", \"60\u00f74=15, 0\"],\n  [\"61\u00f75=12, 1\", \"17\u00f77=2, 3\"],\n  [\"51\u00f76=8, 3\", \"70\u00f72=35, 0\"],\n  [\"45\u00f73=15, 0\", \"49\u00f78=6, 1\"],\n  [\"89\u00f77=12, 5\", \"27\u00f72=13, 1\"],\n  [\"75\u00f75=15, 0\", \"60\u00f78=7, 4\"],\n  [\"49\u00f75=9, 4\", \"58\u00f76=9, 4\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit-division answer cell's text with its new value.\n# The mapping below was derived from the target diff: each entry is\n# old -> new text for a table-cell run in the document.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('42\u00f76=7, 0', '40\u00f72=20, 0'),\n    @('87\u00f78=10, 7', '34\u00f76=5, 4'),\n    @('13\u00f77=1, 6', '47\u00f73=15, 2'),\n    @('60\u00f79=6, 6', '53\u00f73=17, 2'),\n    @('17\u00f78=2, 1', '84\u00f74=21, 0'),\n    @('31\u00f72=15, 1', '32\u00f74=8, 0'),\n    @('29\u00f73=9, 2', '85\u00f73=28, 1'),\n    @('99\u00f78=12, 3', '65\u00f74=16, 1'),\n    @('54\u00f78=6, 6', '42\u00f73=14, 0'),\n    @('31\u00f76=5, 1', '46\u00f79=5, 1'),\n    @('92\u00f74=23, 0', '71\u00f79=7, 8'),\n    @('13\u00f79=1, 4', '45\u00f78=5, 5'),\n    @('83\u00f79=9, 2', '94\u00f76=15, 4'),\n    @('83\u00f73=27, 2', '11\u00f79=1, 2'),\n    @('93\u00f76=15, 3', '20\u00f75=4, 0'),\n    @('90\u00f72=45, 0', '69\u00f78=8, 5'),\n    @('12\u00f77=1, 5', '98\u00f72=49, 0'),\n    @('16\u00f77=2, 2', '68\u00f77=9, 5'),\n    @('24\u00f74=6, 0', '60\u00f74=15, 0'),\n    @('61\u00f75=12, 1', '17\u00f77=2, 3'),\n    @('51\u00f76=8, 3', '70\u00f72=35, 0'),\n    @('45\u00f73=15, 0', '49\u00f78=6, 1'),\n    @('89\u00f77=12, 5', '27\u00f72=13, 1'),\n    @('75\u00f75=15, 0', '60\u00f78=7, 4'),\n    @('49\u00f75=9, 4', '58\u00f76=9, 4'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
